# Replace the 100 arithmetic-problem cell texts with the regenerated
# problem set. MatchWholeWord ($true, 3rd arg) is essential: several of
# the old strings (e.g. "6+75=") are literal substrings of other old
# strings (e.g. "16+75="), and whole-word matching (treating '+'/'='/digits
# boundaries correctly) keeps each replacement scoped to its own cell.
$d = $word.ActiveDocument
$d.Content.Find.Execute("54+3=", $false, $true, $false, $false, $false, $true, 1, $false, "27+15=", 2) | Out-Null
$d.Content.Find.Execute("38-11=", $false, $true, $false, $false, $false, $true, 1, $false, "10-9=", 2) | Out-Null
$d.Content.Find.Execute("94-72=", $false, $true, $false, $false, $false, $true, 1, $false, "68-5=", 2) | Out-Null
$d.Content.Find.Execute("10+69=", $false, $true, $false, $false, $false, $true, 1, $false, "29+10=", 2) | Out-Null
$d.Content.Find.Execute("45-7=", $false, $true, $false, $false, $false, $true, 1, $false, "72-34=", 2) | Out-Null
$d.Content.Find.Execute("79+5=", $false, $true, $false, $false, $false, $true, 1, $false, "67-30=", 2) | Out-Null
$d.Content.Find.Execute("91-28=", $false, $true, $false, $false, $false, $true, 1, $false, "90-59=", 2) | Out-Null
$d.Content.Find.Execute("8+76=", $false, $true, $false, $false, $false, $true, 1, $false, "11+19=", 2) | Out-Null
$d.Content.Find.Execute("86-52=", $false, $true, $false, $false, $false, $true, 1, $false, "38+46=", 2) | Out-Null
$d.Content.Find.Execute("21+24=", $false, $true, $false, $false, $false, $true, 1, $false, "67-30=", 2) | Out-Null
$d.Content.Find.Execute("74-72=", $false, $true, $false, $false, $false, $true, 1, $false, "83-21=", 2) | Out-Null
$d.Content.Find.Execute("19+62=", $false, $true, $false, $false, $false, $true, 1, $false, "34-30=", 2) | Out-Null
$d.Content.Find.Execute("63+32=", $false, $true, $false, $false, $false, $true, 1, $false, "55+33=", 2) | Out-Null
$d.Content.Find.Execute("23+71=", $false, $true, $false, $false, $false, $true, 1, $false, "91-44=", 2) | Out-Null
$d.Content.Find.Execute("9+57=", $false, $true, $false, $false, $false, $true, 1, $false, "3+50=", 2) | Out-Null
$d.Content.Find.Execute("44+14=", $false, $true, $false, $false, $false, $true, 1, $false, "80-69=", 2) | Out-Null
$d.Content.Find.Execute("61+35=", $false, $true, $false, $false, $false, $true, 1, $false, "76-29=", 2) | Out-Null
$d.Content.Find.Execute("45+48=", $false, $true, $false, $false, $false, $true, 1, $false, "9+68=", 2) | Out-Null
$d.Content.Find.Execute("23+14=", $false, $true, $false, $false, $false, $true, 1, $false, "3+15=", 2) | Out-Null
$d.Content.Find.Execute("6+75=", $false, $true, $false, $false, $false, $true, 1, $false, "6+13=", 2) | Out-Null
$d.Content.Find.Execute("72-58=", $false, $true, $false, $false, $false, $true, 1, $false, "38+45=", 2) | Out-Null
$d.Content.Find.Execute("32-8=", $false, $true, $false, $false, $false, $true, 1, $false, "76-45=", 2) | Out-Null
$d.Content.Find.Execute("44+44=", $false, $true, $false, $false, $false, $true, 1, $false, "34+50=", 2) | Out-Null
$d.Content.Find.Execute("14+80=", $false, $true, $false, $false, $false, $true, 1, $false, "82-25=", 2) | Out-Null
$d.Content.Find.Execute("16+75=", $false, $true, $false, $false, $false, $true, 1, $false, "62+14=", 2) | Out-Null
$d.Content.Find.Execute("97-1=", $false, $true, $false, $false, $false, $true, 1, $false, "35+4=", 2) | Out-Null
$d.Content.Find.Execute("75-15=", $false, $true, $false, $false, $false, $true, 1, $false, "61-8=", 2) | Out-Null
$d.Content.Find.Execute("97-92=", $false, $true, $false, $false, $false, $true, 1, $false, "53+22=", 2) | Out-Null
$d.Content.Find.Execute("23-1=", $false, $true, $false, $false, $false, $true, 1, $false, "65+12=", 2) | Out-Null
$d.Content.Find.Execute("0+22=", $false, $true, $false, $false, $false, $true, 1, $false, "10+68=", 2) | Out-Null
$d.Content.Find.Execute("82-70=", $false, $true, $false, $false, $false, $true, 1, $false, "25+44=", 2) | Out-Null
$d.Content.Find.Execute("76-53=", $false, $true, $false, $false, $false, $true, 1, $false, "29-11=", 2) | Out-Null
$d.Content.Find.Execute("36+39=", $false, $true, $false, $false, $false, $true, 1, $false, "20+49=", 2) | Out-Null
$d.Content.Find.Execute("12+10=", $false, $true, $false, $false, $false, $true, 1, $false, "53-4=", 2) | Out-Null
$d.Content.Find.Execute("36+17=", $false, $true, $false, $false, $false, $true, 1, $false, "40-36=", 2) | Out-Null
$d.Content.Find.Execute("39+59=", $false, $true, $false, $false, $false, $true, 1, $false, "28-11=", 2) | Out-Null
$d.Content.Find.Execute("46+34=", $false, $true, $false, $false, $false, $true, 1, $false, "2+41=", 2) | Out-Null
$d.Content.Find.Execute("4+38=", $false, $true, $false, $false, $false, $true, 1, $false, "36+29=", 2) | Out-Null
$d.Content.Find.Execute("91-20=", $false, $true, $false, $false, $false, $true, 1, $false, "50+31=", 2) | Out-Null
$d.Content.Find.Execute("72+25=", $false, $true, $false, $false, $false, $true, 1, $false, "41+11=", 2) | Out-Null
$d.Content.Find.Execute("76-35=", $false, $true, $false, $false, $false, $true, 1, $false, "93-53=", 2) | Out-Null
$d.Content.Find.Execute("97-68=", $false, $true, $false, $false, $false, $true, 1, $false, "97-2=", 2) | Out-Null
$d.Content.Find.Execute("13+65=", $false, $true, $false, $false, $false, $true, 1, $false, "22-4=", 2) | Out-Null
$d.Content.Find.Execute("23+22=", $false, $true, $false, $false, $false, $true, 1, $false, "91-6=", 2) | Out-Null
$d.Content.Find.Execute("93-55=", $false, $true, $false, $false, $false, $true, 1, $false, "71-46=", 2) | Out-Null
$d.Content.Find.Execute("77+14=", $false, $true, $false, $false, $false, $true, 1, $false, "21+40=", 2) | Out-Null
$d.Content.Find.Execute("99-46=", $false, $true, $false, $false, $false, $true, 1, $false, "61-0=", 2) | Out-Null
$d.Content.Find.Execute("98-72=", $false, $true, $false, $false, $false, $true, 1, $false, "25+56=", 2) | Out-Null
$d.Content.Find.Execute("22-11=", $false, $true, $false, $false, $false, $true, 1, $false, "77-0=", 2) | Out-Null
$d.Content.Find.Execute("35-14=", $false, $true, $false, $false, $false, $true, 1, $false, "35+45=", 2) | Out-Null
$d.Content.Find.Execute("19-1=", $false, $true, $false, $false, $false, $true, 1, $false, "34+62=", 2) | Out-Null
$d.Content.Find.Execute("48-34=", $false, $true, $false, $false, $false, $true, 1, $false, "75-54=", 2) | Out-Null
$d.Content.Find.Execute("88-48=", $false, $true, $false, $false, $false, $true, 1, $false, "97-26=", 2) | Out-Null
$d.Content.Find.Execute("4+48=", $false, $true, $false, $false, $false, $true, 1, $false, "15+47=", 2) | Out-Null
$d.Content.Find.Execute("29-25=", $false, $true, $false, $false, $false, $true, 1, $false, "35+7=", 2) | Out-Null
$d.Content.Find.Execute("2+59=", $false, $true, $false, $false, $false, $true, 1, $false, "54-43=", 2) | Out-Null
$d.Content.Find.Execute("20+1=", $false, $true, $false, $false, $false, $true, 1, $false, "37+15=", 2) | Out-Null
$d.Content.Find.Execute("5-2=", $false, $true, $false, $false, $false, $true, 1, $false, "33-16=", 2) | Out-Null
$d.Content.Find.Execute("72-36=", $false, $true, $false, $false, $false, $true, 1, $false, "88-85=", 2) | Out-Null
$d.Content.Find.Execute("23-2=", $false, $true, $false, $false, $false, $true, 1, $false, "17+79=", 2) | Out-Null
$d.Content.Find.Execute("69-48=", $false, $true, $false, $false, $false, $true, 1, $false, "74+18=", 2) | Out-Null
$d.Content.Find.Execute("55-1=", $false, $true, $false, $false, $false, $true, 1, $false, "56-34=", 2) | Out-Null
$d.Content.Find.Execute("13+63=", $false, $true, $false, $false, $false, $true, 1, $false, "1+6=", 2) | Out-Null
$d.Content.Find.Execute("49-3=", $false, $true, $false, $false, $false, $true, 1, $false, "76-59=", 2) | Out-Null
$d.Content.Find.Execute("85-44=", $false, $true, $false, $false, $false, $true, 1, $false, "58+13=", 2) | Out-Null
$d.Content.Find.Execute("10+53=", $false, $true, $false, $false, $false, $true, 1, $false, "19+56=", 2) | Out-Null
$d.Content.Find.Execute("16+49=", $false, $true, $false, $false, $false, $true, 1, $false, "76-39=", 2) | Out-Null
$d.Content.Find.Execute("48-17=", $false, $true, $false, $false, $false, $true, 1, $false, "90-40=", 2) | Out-Null
$d.Content.Find.Execute("43+13=", $false, $true, $false, $false, $false, $true, 1, $false, "40+47=", 2) | Out-Null
$d.Content.Find.Execute("68-14=", $false, $true, $false, $false, $false, $true, 1, $false, "75-29=", 2) | Out-Null
$d.Content.Find.Execute("41+53=", $false, $true, $false, $false, $false, $true, 1, $false, "62-5=", 2) | Out-Null
$d.Content.Find.Execute("76+4=", $false, $true, $false, $false, $false, $true, 1, $false, "73+15=", 2) | Out-Null
$d.Content.Find.Execute("94+2=", $false, $true, $false, $false, $false, $true, 1, $false, "71-57=", 2) | Out-Null
$d.Content.Find.Execute("93-39=", $false, $true, $false, $false, $false, $true, 1, $false, "47+32=", 2) | Out-Null
$d.Content.Find.Execute("77+1=", $false, $true, $false, $false, $false, $true, 1, $false, "11-2=", 2) | Out-Null
$d.Content.Find.Execute("16+27=", $false, $true, $false, $false, $false, $true, 1, $false, "32+48=", 2) | Out-Null
$d.Content.Find.Execute("79+20=", $false, $true, $false, $false, $false, $true, 1, $false, "88+1=", 2) | Out-Null
$d.Content.Find.Execute("0+46=", $false, $true, $false, $false, $false, $true, 1, $false, "1+82=", 2) | Out-Null
$d.Content.Find.Execute("1+12=", $false, $true, $false, $false, $false, $true, 1, $false, "63+31=", 2) | Out-Null
$d.Content.Find.Execute("26-19=", $false, $true, $false, $false, $false, $true, 1, $false, "82-81=", 2) | Out-Null
$d.Content.Find.Execute("68-60=", $false, $true, $false, $false, $false, $true, 1, $false, "63-27=", 2) | Out-Null
$d.Content.Find.Execute("54-41=", $false, $true, $false, $false, $false, $true, 1, $false, "30+33=", 2) | Out-Null
$d.Content.Find.Execute("10+48=", $false, $true, $false, $false, $false, $true, 1, $false, "20+49=", 2) | Out-Null
$d.Content.Find.Execute("26+57=", $false, $true, $false, $false, $false, $true, 1, $false, "84-55=", 2) | Out-Null
$d.Content.Find.Execute("9+12=", $false, $true, $false, $false, $false, $true, 1, $false, "30-1=", 2) | Out-Null
$d.Content.Find.Execute("72+8=", $false, $true, $false, $false, $false, $true, 1, $false, "29+41=", 2) | Out-Null
$d.Content.Find.Execute("52-13=", $false, $true, $false, $false, $false, $true, 1, $false, "84-12=", 2) | Out-Null
$d.Content.Find.Execute("22+40=", $false, $true, $false, $false, $false, $true, 1, $false, "98-63=", 2) | Out-Null
$d.Content.Find.Execute("94-45=", $false, $true, $false, $false, $false, $true, 1, $false, "44+18=", 2) | Out-Null
$d.Content.Find.Execute("16+60=", $false, $true, $false, $false, $false, $true, 1, $false, "63-47=", 2) | Out-Null
$d.Content.Find.Execute("50-10=", $false, $true, $false, $false, $false, $true, 1, $false, "92-13=", 2) | Out-Null
$d.Content.Find.Execute("80-71=", $false, $true, $false, $false, $false, $true, 1, $false, "15+77=", 2) | Out-Null
$d.Content.Find.Execute("71-13=", $false, $true, $false, $false, $false, $true, 1, $false, "93-31=", 2) | Out-Null
$d.Content.Find.Execute("33-10=", $false, $true, $false, $false, $false, $true, 1, $false, "39-29=", 2) | Out-Null
$d.Content.Find.Execute("62-23=", $false, $true, $false, $false, $false, $true, 1, $false, "75-11=", 2) | Out-Null
$d.Content.Find.Execute("56-35=", $false, $true, $false, $false, $false, $true, 1, $false, "9+82=", 2) | Out-Null
$d.Content.Find.Execute("47-42=", $false, $true, $false, $false, $false, $true, 1, $false, "30+42=", 2) | Out-Null
$d.Content.Find.Execute("54-28=", $false, $true, $false, $false, $false, $true, 1, $false, "64-1=", 2) | Out-Null
$d.Content.Find.Execute("11+84=", $false, $true, $false, $false, $false, $true, 1, $false, "55+35=", 2) | Out-Null
$d.Content.Find.Execute("80-61=", $false, $true, $false, $false, $false, $true, 1, $false, "26+18=", 2) | Out-Null
